# ============================================================================
# Edit script: updates the "APISuite" (sheet2) worksheet of TestData.xlsx
# to add new booking-suite test case rows (negative scenarios), adjust
# formatting/fills on the existing rows, extend the table with two blank
# rows, and tweak column/row sizing + selection.
# ============================================================================

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("UISuite")
$ws2 = $wb.Worksheets.Item("APISuite")

# ----------------------------------------------------------------------
# Alignment / theme-color constants (as read back from this COM runtime)
# ----------------------------------------------------------------------
$xlLeft   = -4131
$xlCenter = -4108
$xlVTop    = -4160
$xlVCenter = -4108
$xlVBottom = -4107

# theme color ids (Interior.ThemeColor) -> resulting <fgColor theme="N"/>
#   2 -> theme="0" (Background1/white)
#   8 -> theme="7" (Accent4 / gold)  - no tint
#   9 -> theme="8" (Accent5 / blue)  - no tint

# ----------------------------------------------------------------------
# 1. Update the JSON request body text already referenced by D3
#    (same shared string slot, new formatting / no "bookingid" wrapper)
# ----------------------------------------------------------------------
$jsonPositive = @"
{
    "firstname" : "Jim",
    "lastname" : "Brown",
    "totalprice" : 111,
    "depositpaid" : true,
    "bookingdates" : {
        "checkin" : "2018-01-01",
        "checkout" : "2019-01-01"
    },
    "additionalneeds" : "Breakfast"
}
"@
$ws2.Range("D3").Value2 = $jsonPositive

$jsonNullCheck = @"
{
    "firstname" : null,
    "lastname" : null,
    "totalprice" : null,
    "depositpaid" : null,
    "bookingdates" : {
        "checkin" : null,
        "checkout" : null
    },
    "additionalneeds" : null
}
"@

$jsonInvalidDate = @"
{
    "firstname" : "Jim",
    "lastname" : "Brown",
    "totalprice" : 12,
    "depositpaid" : false,
    "bookingdates" : {
        "checkin" : "2021-13-31",
        "checkout" : "2021-13-32"
    },
    "additionalneeds" : null
}
"@

# ----------------------------------------------------------------------
# 2. Row 2 & Row 3 - restyle A/B/C/D to the "fillId4 / left+vcenter" look
#    Row2 uses the plain sheet1-style (fillId4 center / left) -- copy fmt
#    Row3 keeps its existing content/styles, only row height changes
# ----------------------------------------------------------------------
$ws1.Range("A8").Copy($ws2.Range("A2"))
$ws2.Range("A2").Value2 = 1
$ws1.Range("B8").Copy($ws2.Range("B2"))
$ws2.Range("B2").Value2 = "Test Case 1 - Create Booking Postive scenario"
$ws1.Range("B8").Copy($ws2.Range("C2"))
$ws2.Range("C2").Value2 = "uri"
$ws1.Range("B8").Copy($ws2.Range("D2"))
$ws2.Range("D2").Value2 = "/booking"

$ws1.Range("A8").Copy($ws2.Range("A3"))
$ws2.Range("A3").Value2 = 2
$ws2.Rows.Item(3).RowHeight = 141.75

# ----------------------------------------------------------------------
# 3. Row 4 & Row 5 - continuation of Test Case 1 (Expected Response Code /
#    Booking ID Created), fillId4 look
# ----------------------------------------------------------------------
$ws1.Range("A8").Copy($ws2.Range("A4"))
$ws2.Range("A4").Value2 = 3
$ws2.Range("B3").Copy($ws2.Range("B4"))
$ws2.Range("B4").Value2 = "Test Case 1 - Create Booking Postive scenario"
$ws1.Range("B8").Copy($ws2.Range("C4"))
$ws2.Range("C4").Value2 = "Expected Response Code"
$ws1.Range("B8").Copy($ws2.Range("D4"))
$ws2.Range("D4").Value2 = 200

$ws1.Range("A8").Copy($ws2.Range("A5"))
$ws2.Range("A5").Value2 = 4
$ws2.Range("B3").Copy($ws2.Range("B5"))
$ws2.Range("B5").Value2 = "Test Case 1 - Create Booking Postive scenario"
$ws1.Range("B8").Copy($ws2.Range("C5"))
$ws2.Range("C5").Value2 = "Booking ID Created"
$ws1.Range("B8").Copy($ws2.Range("D5"))
$ws2.Range("D5").Value2 = ""

# ----------------------------------------------------------------------
# 4. Row 6,7,8 - Test Case 2 (Null check), fillId3 look
# ----------------------------------------------------------------------
$ws1.Range("A2").Copy($ws2.Range("A6"))
$ws2.Range("A6").Value2 = 5
$ws1.Range("B2").Copy($ws2.Range("B6"))
$ws2.Range("B6").Value2 = "Test Case 2 - Create Booking Negative scenario Null Check"
$ws2.Range("B6").VerticalAlignment = $xlVCenter
$ws1.Range("B2").Copy($ws2.Range("C6"))
$ws2.Range("C6").Value2 = "uri"
$ws1.Range("B2").Copy($ws2.Range("D6"))
$ws2.Range("D6").Value2 = "/booking"

$ws1.Range("A2").Copy($ws2.Range("A7"))
$ws2.Range("A7").Value2 = 6
$ws2.Range("B6").Copy($ws2.Range("B7"))
$ws2.Range("B7").Value2 = "Test Case 2 - Create Booking Negative scenario Null Check"
$ws2.Range("B6").Copy($ws2.Range("C7"))
$ws2.Range("C7").Value2 = "Request Body"
$ws2.Range("D3").Copy($ws2.Range("D7"))
$ws2.Range("D7").Value2 = $jsonNullCheck
$ws2.Range("D7").Interior.ThemeColor = 8
$ws2.Range("D7").VerticalAlignment = $xlVBottom
$ws2.Rows.Item(7).RowHeight = 141.75

$ws1.Range("A2").Copy($ws2.Range("A8"))
$ws2.Range("A8").Value2 = 7
$ws1.Range("B2").Copy($ws2.Range("B8"))
$ws2.Range("B8").Value2 = "Test Case 2 - Create Booking Negative scenario Null Check"
$ws2.Range("B8").VerticalAlignment = $xlVCenter
$ws1.Range("B2").Copy($ws2.Range("C8"))
$ws2.Range("C8").Value2 = "Expected Response Code"
$ws1.Range("B2").Copy($ws2.Range("D8"))
$ws2.Range("D8").Value2 = 500

# ----------------------------------------------------------------------
# 5. Row 9,10,11,12 - Test Case 3 (Invalid date check), fillId4 look
# ----------------------------------------------------------------------
$ws1.Range("A8").Copy($ws2.Range("A9"))
$ws2.Range("A9").Value2 = 8
$ws2.Range("B4").Copy($ws2.Range("B9"))
$ws2.Range("B9").Value2 = "Test Case 3 - Create Booking Negative scenario Invalid Date Check"
$ws2.Range("C4").Copy($ws2.Range("C9"))
$ws2.Range("C9").Value2 = "uri"
$ws2.Range("C4").Copy($ws2.Range("D9"))
$ws2.Range("D9").Value2 = "/booking"

$ws1.Range("A8").Copy($ws2.Range("A10"))
$ws2.Range("A10").Value2 = 9
$ws2.Range("B9").Copy($ws2.Range("B10"))
$ws2.Range("B10").Value2 = "Test Case 3 - Create Booking Negative scenario Invalid Date Check"
$ws2.Range("B9").Copy($ws2.Range("C10"))
$ws2.Range("C10").Value2 = "Request Body"
$ws2.Range("D3").Copy($ws2.Range("D10"))
$ws2.Range("D10").Value2 = $jsonInvalidDate
$ws2.Range("D10").Interior.ThemeColor = 9
$ws2.Range("D10").VerticalAlignment = $xlVCenter
$ws2.Rows.Item(10).RowHeight = 140.25

$ws1.Range("A8").Copy($ws2.Range("A11"))
$ws2.Range("A11").Value2 = 10
$ws2.Range("B9").Copy($ws2.Range("B11"))
$ws2.Range("B11").Value2 = "Test Case 3 - Create Booking Negative scenario Invalid Date Check"
$ws2.Range("B9").Copy($ws2.Range("C11"))
$ws2.Range("C11").Value2 = "Expected Response Code"
$ws2.Range("B9").Copy($ws2.Range("D11"))
$ws2.Range("D11").Value2 = 200

$ws1.Range("A8").Copy($ws2.Range("A12"))
$ws2.Range("A12").Value2 = 11
$ws2.Range("B9").Copy($ws2.Range("B12"))
$ws2.Range("B12").Value2 = "Test Case 3 - Create Booking Negative scenario Invalid Date Check"
$ws2.Range("C4").Copy($ws2.Range("C12"))
$ws2.Range("C12").Value2 = "Expected Response Message"
$ws2.Range("C4").Copy($ws2.Range("D12"))
$ws2.Range("D12").Value2 = "Invalid date"

# ----------------------------------------------------------------------
# 6. Rows 13-17 - renumber the trailing blank rows (A column only),
#    switching their style to the white-filled "fillId5" center look
# ----------------------------------------------------------------------
$ws2.Range("A13").Interior.ThemeColor = 2
$ws2.Range("A14").Interior.ThemeColor = 2
$ws2.Range("A15").Interior.ThemeColor = 2
$ws2.Range("A16").Interior.ThemeColor = 2
$ws2.Range("A16").Value2 = 15
$ws2.Range("A17").Interior.ThemeColor = 2
$ws2.Range("A17").Value2 = 16

# ----------------------------------------------------------------------
# 7. Extend the sheet with two new blank rows (31 & 32), same blank
#    style as rows 18-30
# ----------------------------------------------------------------------
$ws2.Range("A30:D30").Copy($ws2.Range("A31:D31"))
$ws2.Range("A31").Value2 = ""
$ws2.Range("A30:D30").Copy($ws2.Range("A32:D32"))
$ws2.Range("A32").Value2 = ""

# ----------------------------------------------------------------------
# 8. Column B width & selection
# ----------------------------------------------------------------------
$ws2.Columns.Item(2).ColumnWidth = 74.25
$ws2.Range("B10").Select()

Write-Host "Edit complete"
